{"js": "// The cover-letter author list originally read \"...Ryan R, Brinkman...\"\n// (a stray comma instead of the intended abbreviation-period after the\n// middle initial \"R\"). The final revision fixes this to\n// \"...Ryan R. Brinkman...\" while leaving every other word and all\n// character formatting (font, size, language, etc.) untouched.\n//\n// Search for the exact, narrowly-scoped substring \"Ryan R,\" inside the\n// paragraph that lists the authors, then replace just that match with\n// \"Ryan R.\" so the rest of the sentence is left completely intact.\nconst searchResults = context.document.body.search(\"Ryan R,\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find \"Ryan R,\" in the document body.');\n}\n\n// Replace the matched range's text in place; insertText with the\n// \"Replace\" location substitutes the text of the found range while\n// preserving the run's existing character formatting.\nsearchResults.items[0].insertText(\"Ryan R.\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The cover-letter's author list originally read \"...Ryan R, Brinkman...\"\n# (a stray comma where the middle-initial abbreviation period belongs).\n# The final revision corrects this single typo to \"...Ryan R. Brinkman...\",\n# leaving every other word and all character formatting (font, size,\n# language, etc.) untouched.\n#\n# Word find/replace constants (wdReplaceOne / wdFindContinue) spelled out\n# for readability instead of using bare magic numbers.\n$wdFindContinue = 1\n$wdReplaceOne   = 1\n\n$d = $word.ActiveDocument\n\n# Search the whole story (Content range) for the exact, narrowly-scoped\n# text \"Ryan R,\" and replace just that one occurrence with \"Ryan R.\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Ryan R,\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Ryan R.\"\n$find.Forward = $true\n$find.Wrap = $wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$found = $find.Execute(\n    $find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    $wdFindContinue,\n    $false,\n    $find.Replacement.Text,\n    $wdReplaceOne\n)\n\nif (-not $found) {\n    throw 'Could not find \"Ryan R,\" in the document content.'\n}\n"}
